$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: "Parcel (ton)" header plus per-row tonnage values
$ws.Range("G1").Value = "Parcel (ton)"

$parcelValues = @(12,7,11,8,16,5,3,19,10,5,12,16,5,3,19,10,5,12,20,18,3,6,9,3,19,23,25,12,16,5,3,12,7,11,8,5,16,13,16,22,19,13,6)

for ($i = 0; $i -lt $parcelValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $parcelValues[$i]
}

# Match the author's final view state (scrolled/selected cell)
$ws.Range("G57").Select() | Out-Null
